# Updates crypto price (column D) and volume-change (column E) values
# to reflect the latest scrape, per commit "Updated cryptos list on
# Sun May 14 21:25:58 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.750.49"
$ws.Range("E2").Value = "  +0.66%  "
$ws.Range("D3").Value = "1.850.80"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("D4").Value = "'1.034"
$ws.Range("E4").Value = "  +0.35%  "
$ws.Range("D5").Value = "'322.03"
$ws.Range("E5").Value = "  +0.57%  "
$ws.Range("D6").Value = "'1.029"
$ws.Range("E6").Value = "  +0.25%  "
$ws.Range("D7").Value = "'0.4396"
$ws.Range("E7").Value = "  +0.58%  "
$ws.Range("D8").Value = "'0.3807"
$ws.Range("E8").Value = "  +1.85%  "
$ws.Range("D9").Value = "'0.07408"
$ws.Range("E9").Value = "  +0.34%  "
$ws.Range("D10").Value = "'0.8856"
$ws.Range("E10").Value = "  +1.25%  "
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.856.01"
$ws.Range("E12").Value = "  +0.35%  "
$ws.Range("D13").Value = "'5.511"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("D14").Value = "'6.698"
$ws.Range("E14").Value = "  +0.43%  "
$ws.Range("D15").Value = "'0.07177"
$ws.Range("E15").Value = "  +0.36%  "
$ws.Range("D16").Value = "'85.15"
$ws.Range("E16").Value = "  +2.89%  "
$ws.Range("D17").Value = "'1.037"
$ws.Range("E17").Value = "  +0.51%  "
$ws.Range("D18").Value = "'0.000009085"
$ws.Range("E18").Value = "  +0.77%  "
$ws.Range("D19").Value = "'1.030"
$ws.Range("E19").Value = "  +0.34%  "
$ws.Range("D20").Value = "'15.50"
$ws.Range("E20").Value = "  +0.65%  "
$ws.Range("D21").Value = "27.740.65"
$ws.Range("E21").Value = "  +0.64%  "
$ws.Range("D22").Value = "'5.269"
$ws.Range("E22").Value = "  +0.54%  "
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "2.082.21"
$ws.Range("E24").Value = "  +0.65%  "
$ws.Range("D25").Value = "'2.057"
$ws.Range("E25").Value = "  +6.53%  "
$ws.Range("D26").Value = "'158.26"
$ws.Range("E26").Value = "  +0.68%  "
$ws.Range("D27").Value = "'18.70"
$ws.Range("E27").Value = "  +0.01%  "
$ws.Range("D28").Value = "'1.996"
$ws.Range("D29").Value = "'5.337"
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").Value = "'118.11"
$ws.Range("E30").Value = "  +1.81%  "
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("E32").Value = "  +0.57%  "
$ws.Range("D33").Value = "'0.7730"
$ws.Range("E33").Value = "  +0.99%  "
$ws.Range("D34").Value = "'3.001"
$ws.Range("E34").Value = "  +4.32%  "
$ws.Range("D35").Value = "'4.573"
$ws.Range("E35").Value = "  +1.54%  "
$ws.Range("E36").Value = "  +0.31%  "
$ws.Range("D37").Value = "'1.150"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("D38").Value = "'0.01973"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'0.05288"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Value = "'2.861"
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").Value = "'0.5177"
$ws.Range("E41").Value = "  +0.29%  "
$ws.Range("D42").Value = "'0.1670"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'6.865"
$ws.Range("E43").Value = "  +2.77%  "
$ws.Range("D44").Value = "'8.751"
$ws.Range("E44").Value = "  +2.38%  "
$ws.Range("D45").Value = "'110.54"
$ws.Range("E45").Value = "  +1.65%  "
$ws.Range("D46").Value = "'10.77"
$ws.Range("E46").Value = "  +1.76%  "
$ws.Range("E47").Value = "  +0.36%  "
$ws.Range("D48").Value = "'0.06577"
$ws.Range("E48").Value = "  +2.98%  "
$ws.Range("D49").Value = "'1.708"
$ws.Range("E49").Value = "  +0.00%  "
$ws.Range("D50").Value = "'0.4709"
$ws.Range("E50").Value = "  +1.42%  "
$ws.Range("D51").Value = "'1.892"
$ws.Range("E51").Value = "  +0.63%  "
